# Apply cryptos list update (Wed Nov 15 14:40:27 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.124.63"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "2.012.26"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'252.18"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "'0.642"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("D7").Value = "'61.95"
$ws.Range("E7").Value = "  +8.25%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'59.04"
$ws.Range("E9").Value = "  -6.36%  "
$ws.Range("D10").Value = "'0.370"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'0.0746"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'0.910"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'14.82"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "2.309.96"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "'5.41"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "'19.44"
$ws.Range("E17").Value = "  +8.93%  "
$ws.Range("D18").Value = "2.013.03"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "36.082.42"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "'71.94"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'5.26"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").Value = "'233.87"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +17.55%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("D27").Value = "'9.48"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'164.18"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'19.56"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'5.09"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.20"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  +25.03%  "
$ws.Range("D34").Value = "'0.0602"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("E35").Value = "  +12.72%  "
$ws.Range("D36").Value = "'4.49"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "'5.81"
$ws.Range("E39").Value = "  +13.79%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "  +14.36%  "
$ws.Range("D41").Value = "'1.21"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'2.89"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "'0.0215"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "'1.12"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "'16.57"
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("D46").Value = "'7.82"
$ws.Range("E46").Value = "  +4.88%  "
$ws.Range("D47").Value = "'93.52"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "1.417.87"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "'2.48"
$ws.Range("E49").Value = "  +9.63%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "'47.33"
$ws.Range("E51").Value = "  +2.98%  "
